$wb = $excel.ActiveWorkbook
$wsPrice = $wb.Worksheets.Item("plumPriceVersion")
$wsFitting = $wb.Worksheets.Item("FittingPriceInfo")
$wsPrice.Move($wsFitting)
